$d = $word.ActiveDocument

# 1) "This project, developed by a group of 12 undergraduate, was created to
#    simplify this process." -> insert "students" between "undergraduate"
#    and ", was" (so the sentence reads "... undergraduate students, was ...")
$r1 = $d.Content
$r1.Find.Execute("undergraduate, was", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "undergraduate students, was", 2)

# 2) "... It takes specifications and the classes ..." -> "preferences"
$r2 = $d.Content
$r2.Find.Execute("specifications", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "preferences", 2)

# 3) "Specifications can be made by the student ..." -> "Preferences"
$r3 = $d.Content
$r3.Find.Execute("Specifications", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "Preferences", 2)

# 4) "... if a certain specification made makes a sequence ..." -> "preferences"
$r4 = $d.Content
$r4.Find.Execute("specification made", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "preferences made", 2)
